$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5735.6665
$ws.Range("I33").Value = 86.2
$ws.Range("J33").Value = 12797.5
$ws.Range("K33").Value = 86.2
$ws.Range("L33").Value = 12797.5
$ws.Range("M33").Value = 142.8
$ws.Range("N33").Value = -13255.5
$ws.Range("H64").Value = 4379.1
$ws.Range("I64").Value = 4213
$ws.Range("J64").Value = 4489.8335
$ws.Range("K64").Value = 4213
$ws.Range("L64").Value = 4489.8335
$ws.Range("M64").Value = -3965
$ws.Range("N64").Value = -4985.8335
$ws.Range("H67").Value = 4379.1
$ws.Range("I67").Value = 4213
$ws.Range("J67").Value = 4489.8335
$ws.Range("K67").Value = 4213
$ws.Range("L67").Value = 4489.8335
$ws.Range("M67").Value = -3355
$ws.Range("N67").Value = -6205.8335
$ws.Range("H112").Value = 5927.485
$ws.Range("J112").Value = 6871.6787
$ws.Range("L112").Value = 20615.0361
$ws.Range("N112").Value = -22831.0361
$ws.Range("H132").Value = 1193.5106
$ws.Range("I132").Value = 1122.6511
$ws.Range("J132").Value = 1955.25
$ws.Range("K132").Value = 3367.9533
$ws.Range("L132").Value = 5865.75
$ws.Range("M132").Value = -837.9533000000001
$ws.Range("N132").Value = -10925.75
$ws.Range("H137").Value = 1532.28
$ws.Range("I137").Value = 1013.86664
$ws.Range("J137").Value = 2309.9
$ws.Range("K137").Value = 3041.59992
$ws.Range("L137").Value = 6929.700000000001
$ws.Range("M137").Value = -491.5999199999997
$ws.Range("N137").Value = -12029.7
$ws.Range("H138").Value = 4199.343
$ws.Range("I138").Value = 2737.12
$ws.Range("J138").Value = 5011.689
$ws.Range("K138").Value = 8211.360000000001
$ws.Range("L138").Value = 15035.067
$ws.Range("M138").Value = -3071.360000000001
$ws.Range("N138").Value = -25315.067

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 977161.6
$ws.Range("I32").Value = 10571.762
$ws.Range("J32").Value = 14509420
$ws.Range("K32").Value = 10571.762
$ws.Range("L32").Value = 14509420
$ws.Range("M32").Value = -10284.762
$ws.Range("N32").Value = -14509994
$ws.Range("H45").Value = 2075.2654
$ws.Range("I45").Value = 1979.6285
$ws.Range("J45").Value = 2314.3572
$ws.Range("K45").Value = 1979.6285
$ws.Range("L45").Value = 2314.3572
$ws.Range("M45").Value = -1602.6285
$ws.Range("N45").Value = -3068.3572
$ws.Range("H102").Value = 2730
$ws.Range("I102").Value = 1600
$ws.Range("J102").Value = 7250
$ws.Range("K102").Value = 1600
$ws.Range("L102").Value = 7250
$ws.Range("M102").Value = 22
$ws.Range("N102").Value = -10494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1660
$ws.Range("I105").Value = 1660
$ws.Range("K105").Value = 1660
$ws.Range("M105").Value = 87
$ws.Range("H134").Value = 912.625
$ws.Range("I134").Value = 870.3333
$ws.Range("J134").Value = 1039.5
$ws.Range("K134").Value = 2610.9999
$ws.Range("L134").Value = 3118.5
$ws.Range("M134").Value = -75.9998999999998
$ws.Range("N134").Value = -8188.5
$ws.Range("H140").Value = 86740
$ws.Range("J140").Value = 86740
$ws.Range("L140").Value = 86740
$ws.Range("N140").Value = -97100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 645.125
$ws.Range("I2").Value = 594.5714
$ws.Range("J2").Value = 999
$ws.Range("K2").Value = 594.5714
$ws.Range("L2").Value = 999
$ws.Range("M2").Value = -481.5714
$ws.Range("N2").Value = -1225
$ws.Range("H62").Value = 3882.7144
$ws.Range("I62").Value = 3700.625
$ws.Range("J62").Value = 4125.5
$ws.Range("K62").Value = 3700.625
$ws.Range("L62").Value = 4125.5
$ws.Range("M62").Value = -3076.625
$ws.Range("N62").Value = -5373.5
$ws.Range("H65").Value = 3882.7144
$ws.Range("I65").Value = 3700.625
$ws.Range("J65").Value = 4125.5
$ws.Range("K65").Value = 18503.125
$ws.Range("L65").Value = 20627.5
$ws.Range("M65").Value = -15383.125
$ws.Range("N65").Value = -26867.5
$ws.Range("H141").Value = 50000
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11112298
$ws.Range("J131").Value = 12196123
$ws.Range("L131").Value = 36588369
$ws.Range("N131").Value = -36598449

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 711066.7
$ws.Range("I11").Value = 2650000
$ws.Range("J11").Value = 6000
$ws.Range("K11").Value = 2650000
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = -2649861
$ws.Range("N11").Value = -6278
$ws.Range("H70").Value = 20841736
$ws.Range("I70").Value = 31257328
$ws.Range("J70").Value = 10555.556
$ws.Range("K70").Value = 31257328
$ws.Range("L70").Value = 10555.556
$ws.Range("M70").Value = -31257058
$ws.Range("N70").Value = -11095.556
$ws.Range("H73").Value = 20841736
$ws.Range("I73").Value = 31257328
$ws.Range("J73").Value = 10555.556
$ws.Range("K73").Value = 31257328
$ws.Range("L73").Value = 10555.556
$ws.Range("M73").Value = -31256392
$ws.Range("N73").Value = -12427.556
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H138").Value = 69650
$ws.Range("J138").Value = 69650
$ws.Range("L138").Value = 69650
$ws.Range("N138").Value = -79930
$ws.Range("H141").Value = 68000
$ws.Range("J141").Value = 68000
$ws.Range("L141").Value = 68000
$ws.Range("N141").Value = -78360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1212
$ws.Range("I2").Value = 1015
$ws.Range("K2").Value = 1015
$ws.Range("M2").Value = -903
$ws.Range("H7").Value = 3777.1904
$ws.Range("I7").Value = 3010
$ws.Range("J7").Value = 4474.636
$ws.Range("K7").Value = 3010
$ws.Range("L7").Value = 4474.636
$ws.Range("M7").Value = -2898
$ws.Range("N7").Value = -4698.636
$ws.Range("H68").Value = 2054.9119
$ws.Range("I68").Value = 1963.5
$ws.Range("J68").Value = 2185.5
$ws.Range("K68").Value = 1963.5
$ws.Range("L68").Value = 2185.5
$ws.Range("M68").Value = -1214.5
$ws.Range("N68").Value = -3683.5
$ws.Range("H71").Value = 2054.9119
$ws.Range("I71").Value = 1963.5
$ws.Range("J71").Value = 2185.5
$ws.Range("K71").Value = 9817.5
$ws.Range("L71").Value = 10927.5
$ws.Range("M71").Value = -6073.5
$ws.Range("N71").Value = -18415.5
$ws.Range("H100").Value = 3023.2
$ws.Range("I100").Value = 2475.8462
$ws.Range("J100").Value = 3616.1667
$ws.Range("K100").Value = 2475.8462
$ws.Range("L100").Value = 3616.1667
$ws.Range("M100").Value = -1934.8462
$ws.Range("N100").Value = -4698.1667
$ws.Range("H126").Value = 3777.1904
$ws.Range("I126").Value = 3010
$ws.Range("J126").Value = 4474.636
$ws.Range("K126").Value = 9030
$ws.Range("L126").Value = 13423.908
$ws.Range("M126").Value = -6560
$ws.Range("N126").Value = -18363.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 5000
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 5000
$ws.Range("N39").Value = -5826
